# Update "想去人数" (want-to-go count) figures in sheet "展览" (Exhibitions)
# and sheet "全部类型" (All types) to the freshly scraped values.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15019
$ws1.Range("F3").Value = 19071
$ws1.Range("F5").Value = 139
$ws1.Range("F15").Value = 220
$ws1.Range("F22").Value = 7952
$ws1.Range("F26").Value = 63
$ws1.Range("F27").Value = 1247
$ws1.Range("F28").Value = 19
$ws1.Range("F29").Value = 6060
$ws1.Range("F34").Value = 286
$ws1.Range("F35").Value = 5438
$ws1.Range("F36").Value = 342

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15019
$ws4.Range("F3").Value = 19071
$ws4.Range("F5").Value = 139
$ws4.Range("F15").Value = 220
$ws4.Range("F23").Value = 7952
$ws4.Range("F27").Value = 63
$ws4.Range("F28").Value = 1247
$ws4.Range("F29").Value = 19
$ws4.Range("F32").Value = 6060
$ws4.Range("F37").Value = 286
$ws4.Range("F38").Value = 5438
$ws4.Range("F39").Value = 342
